# Rename the sheet, insert the new "matchNo" column, and add the scraped rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Manan Vohra"

# Insert a new column A for "matchNo" - shifts existing A:L data to B:M
$ws.Range("A1").EntireColumn.Insert()

# Final table data (header + 4 data rows), in display order.
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

$data = @(
    @("16th","Rajasthan Royals","Manan Vohra","c Richardson b Jamieson","7","9","1","0","77.77","Royal Challengers Bangalore","Wankhede","April 22","RCB won by 10 wickets (with 21 balls remaining)"),
    @("4th","Rajasthan Royals","Manan Vohra","c & b Arshdeep Singh","12","8","1","1","150.00","Punjab Kings","Wankhede","April 12","Punjab Kings won by 4 runs"),
    @("7th","Rajasthan Royals","Manan Vohra","c Rabada b Woakes","9","11","2","0","81.81","Delhi Capitals","Wankhede","April 15","Royals won by 3 wickets (with 2 balls remaining)"),
    @("12th","Rajasthan Royals","Manan Vohra","c Jadeja b Curran","14","11","1","1","127.27","Chennai Super Kings","Wankhede","April 19","Super Kings won by 45 runs")
)

# The scraped source stores every cell (even number-looking ones, like "9" or
# "81.81") as plain text. Columns E:I ("runs","balls","fours","sixes","sr")
# are the only ones whose values look fully numeric to Excel's auto-detect,
# so format just those as Text first to stop them turning into real numbers
# (which would also lose the "150.00"-style trailing zeros).
$ws.Range("E1:I5").NumberFormat = "@"

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
